$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# 1) Update existing rows 5-10: content edits (styles already correct)
# ============================================================
$ws.Range('A5').Value = 'TC_1'
$ws.Range('B5').Value = 'Sign in'
$ws.Range('C5').Value = 'To Verify the login page that when  both Sign in & password field is blank and Sign in button is clicked.'
$ws.Range('D5').Value = '1. Open Application                                       2.Leave Username and password blank                                         3.Click sign in Button'
$ws.Range('E5').Value = ''

$ws.Range('A6').Value = 'TC_2'
$ws.Range('B6').Value = 'Sign in'
$ws.Range('C6').Value = 'Verify user able to Sign in with Right Username & Empty Password'
$ws.Range('D6').Value = '1. Open Application                                       2.Leave Username Empty .                                                    3. Enter Password                                                                         4.Click sign in Button'

$ws.Range('A7').Value = 'TC_3'
$ws.Range('B7').Value = 'Sign in'
$ws.Range('C7').Value = 'Verify user able to Sign in with Right Username & Empty Password'
$ws.Range('D7').Value = '1. Open Application                                       2.Enter Username  .                                                    3. Leave Password Empty.                                                                        4.Click sign in Button'

$ws.Range('A8').Value = 'TC_4'
$ws.Range('B8').Value = 'Sign in'
$ws.Range('C8').Value = 'Verify user able to Sign in with Valid Username & Invalid Password'
$ws.Range('D8').Value = '1. Open Application                                       2.Enter valid Username  .                                                    3.Enter invalid Password                                                                      4.Click sign in Button'

$ws.Range('A9').Value = 'TC_5'
$ws.Range('B9').Value = 'Sign in'
$ws.Range('C9').Value = 'Verify user able to Sign in with invalid Username & valid Password'
$ws.Range('D9').Value = '1. Open Application                                       2.Enter Username  .                                                    3. Leave Password Empty.                                                                        4.Click sign in Button'

$ws.Range('A10').Value = 'TC_6'
$ws.Range('B10').Value = 'Sign in'
$ws.Range('C10').Value = 'Verify user  able to Sign in with valid  Username & valid Password'

# C5 picks up the same visual style as C6:C10 (font/wrap) in the edited file
$ws.Range('C6').Copy()
$ws.Range('C5').PasteSpecial(-4122)
$ws.Range('C5').Value = 'To Verify the login page that when  both Sign in & password field is blank and Sign in button is clicked.'

# ============================================================
# 2) Append new test-case rows 11-20 (new content from this upload)
# ============================================================

# Prime formats once: column A/B/D inherit row 10 formatting, column C uses a plain wrapped style
$ws.Range('A10:B10').Copy()
$ws.Range('A11:B20').PasteSpecial(-4122)
$ws.Range('D10').Copy()
$ws.Range('D11:D20').PasteSpecial(-4122)
$ws.Range('C11:C20').WrapText = $true

$ws.Range('A11').Value = 'TC_7 '
$ws.Range('B11').Value = 'Dashboad'
$ws.Range('C11').Value = 'Verify that after login user ab to navigate to Dashboard page'

$ws.Range('A12').Value = 'TC_8'
$ws.Range('B12').Value = 'Dashboad'
$ws.Range('C12').Value = 'verify the functionality of Punch in and Punch out button '

$ws.Range('A13').Value = 'TC_9'
$ws.Range('B13').Value = 'Dashboad'
$ws.Range('C13').Value = 'Verify the functionality for break button '

$ws.Range('A14').Value = 'TC_10'
$ws.Range('B14').Value = 'Dashboad'
$ws.Range('C14').Value = 'Verify the functionality for Add task.'

# Author went back and filled in / corrected the Test Steps (D) column for rows 10-14
$ws.Range('D10').Value = '1. Open Application                                       2.Enter Username  .                                                    3. Enter Password .                                                                        4.Click sign in Button'
$ws.Range('D11').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to dashboard page.'
$ws.Range('D12').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to dashboard page.                  6. click on Punch in button .                                         7. Record Time .                                                       8. Then punch out.'
$ws.Range('D13').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to dashboard page.                  6. click on Punch in button .                                         7. Record Time .                                                       8. Then take a break and verify.'
$ws.Range('D14').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to dashboard page.                  6. Click on Add task and enter any task for the day .'

$ws.Rows(11).RowHeight = 75
$ws.Rows(12).RowHeight = 120
$ws.Rows(13).RowHeight = 120
$ws.Rows(14).RowHeight = 105

$ws.Range('A15').Value = 'TC_11'
$ws.Range('B15').Value = 'Dashboad'
$ws.Range('C15').Value = 'Verify that user is able to add multiple task a day.'
$ws.Range('D15').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to dashboard page.                  6.Add task multple task for the day .'
$ws.Rows(15).RowHeight = 90

$ws.Range('A16').Value = 'TC_12'
$ws.Range('B16').Value = 'Dashboad'
$ws.Range('C16').Value = 'Verify that user ia able to edit task .'
$ws.Range('D16').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to dashboard page.                  6. Click on Add task and enter any task.            7.then click on edit button.'
$ws.Rows(16).RowHeight = 105

$ws.Range('A17').Value = 'TC_13'
$ws.Range('B17').Value = 'Dashboard'
$ws.Range('C17').Value = 'Verify that user is able to Book Meeting room .'
$ws.Range('D17').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to dashboard page.                  6. Click on Book a meeting on meeting manager section'
$ws.Rows(17).RowHeight = 105

$ws.Range('A18').Value = 'TC_14'
$ws.Range('B18').Value = 'Task'
$ws.Range('C18').Value = 'Verify that user is able to see there task on Task page'
$ws.Range('D18').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to Task page.                  '
$ws.Rows(18).RowHeight = 75

$ws.Range('A19').Value = 'TC_15'
$ws.Range('B19').Value = 'Attendance'
$ws.Range('C19').Value = 'Verify that user is able to see ther Attendance on Attendance Page'
$ws.Range('D19').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to Attendance page.                  '
$ws.Rows(19).RowHeight = 75

$ws.Range('A20').Value = 'TC_16'
$ws.Range('B20').Value = 'Dashboad'
$ws.Range('C20').Value = 'Verify that user ia able Logout.'
$ws.Range('D20').Value = '1. Open Application                                       2.Enter Username  .                                                    3.  Enter Password.                                                                        4.Click sign in Button.                                        5.Navigate to Dashbaord page.                            6.Click on Logout.           '
$ws.Rows(20).RowHeight = 90

# ============================================================
# 3) Selection / view state to match the saved file
# ============================================================
$ws.Range('D20').Select()
